$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (B and E get wider to fit the new descriptive text)
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 22.9
$ws.Columns("E").ColumnWidth = 32.45

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 42.5
$ws.Rows(2).RowHeight = 56
$ws.Rows(3).RowHeight = 28
$ws.Rows(4).RowHeight = 28

# ---------------------------------------------------------------------------
# Header row tweaks: switch the header font to "Calibri " and recentre E1
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Font.Name = "Calibri "
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").VerticalAlignment = -4108

$ws.Range("E1").Font.Name = "Calibri "
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").WrapText = $true

# F1/G1 no longer carry any special formatting - clear them out completely
$ws.Range("F1:G1").Clear()

# H1/H2 keep the centred + wrapped look (same visual style as before)
$ws.Range("H1:H2").HorizontalAlignment = -4108
$ws.Range("H1:H2").WrapText = $true

# ---------------------------------------------------------------------------
# Row 2 - PEJABAT NEGARA, PEJABAT DAERAH
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("A2").Font.Name = "Calibri "
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").NumberFormat = "0"
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").ShrinkToFit = $true

$ws.Range("B2").Value = "PEJABAT NEGARA, PEJABAT DAERAH"
$ws.Range("B2").Font.Name = "Calibri "
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").VerticalAlignment = -4160
$ws.Range("B2").WrapText = $true

$ws.Range("C2").Value = "OH"
$ws.Range("C2").Font.Name = "Calibri "
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").VerticalAlignment = -4160
$ws.Range("C2").WrapText = $true

$ws.Range("D2").Value = 250000
$ws.Range("D2").Font.Name = "Calibri "
$ws.Range("D2").Font.Color = 0
$ws.Range("D2").NumberFormat = "#,##0"
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").VerticalAlignment = -4160
$ws.Range("D2").ShrinkToFit = $true

$ws.Range("E2").Value = 125000
$ws.Range("E2").Font.Name = "Calibri "
$ws.Range("E2").Font.Color = 0
$ws.Range("E2").NumberFormat = "#,##0"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4160
$ws.Range("E2").ShrinkToFit = $true

# ---------------------------------------------------------------------------
# Row 3 - PEJABAT ESELON I
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("A3").Font.Name = "Calibri "
$ws.Range("A3").Font.Color = 0
$ws.Range("A3").NumberFormat = "0"
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").ShrinkToFit = $true

$ws.Range("B3").Value = "PEJABAT ESELON  I"
$ws.Range("B3").Font.Name = "Calibri "
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").VerticalAlignment = -4160
$ws.Range("B3").WrapText = $true

$ws.Range("C3").Value = "OH"
$ws.Range("C3").Font.Name = "Calibri "
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4160
$ws.Range("C3").WrapText = $true

$ws.Range("D3").Value = 200000
$ws.Range("D3").Font.Name = "Calibri "
$ws.Range("D3").Font.Color = 0
$ws.Range("D3").NumberFormat = "#,##0"
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4160
$ws.Range("D3").ShrinkToFit = $true

$ws.Range("E3").Value = 100000
$ws.Range("E3").Font.Name = "Calibri "
$ws.Range("E3").Font.Color = 0
$ws.Range("E3").NumberFormat = "#,##0"
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4160
$ws.Range("E3").ShrinkToFit = $true

# ---------------------------------------------------------------------------
# Row 4 - PEJABAT ESELON II
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("A4").Font.Name = "Calibri "
$ws.Range("A4").Font.Color = 0
$ws.Range("A4").NumberFormat = "0"
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("A4").ShrinkToFit = $true

$ws.Range("B4").Value = "PEJABAT ESELON  II"
$ws.Range("B4").Font.Name = "Calibri "
$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("B4").WrapText = $true

$ws.Range("C4").Value = "OH"
$ws.Range("C4").Font.Name = "Calibri "
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true

$ws.Range("D4").Value = 150000
$ws.Range("D4").Font.Name = "Calibri "
$ws.Range("D4").Font.Color = 0
$ws.Range("D4").NumberFormat = "#,##0"
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").VerticalAlignment = -4160
$ws.Range("D4").ShrinkToFit = $true

$ws.Range("E4").Value = 75000
$ws.Range("E4").Font.Name = "Calibri "
$ws.Range("E4").Font.Color = 0
$ws.Range("E4").NumberFormat = "#,##0"
$ws.Range("E4").HorizontalAlignment = -4108
$ws.Range("E4").VerticalAlignment = -4160
$ws.Range("E4").ShrinkToFit = $true

# ---------------------------------------------------------------------------
# Final selection, matching the saved workbook's cursor position
# ---------------------------------------------------------------------------
$ws.Range("B2").Select()

Write-Output "applied"
